# Introduce a separate fee_currency column across the order/transfer sheets,
# consolidate currency_conversions source/target fee columns into a single
# fees + fee_currency pair, and switch the active tab from
# currency_conversions (sheet5) to espp (sheet7).

$wb = $excel.ActiveWorkbook

$xlPasteFormats = -4122

function Set-HeaderStyle($ws, $srcCell, $dstCell) {
    $ws.Range($srcCell).Copy()
    $ws.Range($dstCell).PasteSpecial($xlPasteFormats)
    $excel.CutCopyMode = 0
}

# ---------------------------------------------------------------------------
# buy_orders: date, symbol, quantity, buy_price, fees, currency, comment
#          -> date, symbol, quantity, buy_price, currency, fees, fee_currency, comment
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("buy_orders")

$ws.Cells.Item(1,5).Value = "currency"
$ws.Cells.Item(1,6).Value = "fees"
$ws.Cells.Item(1,7).Value = "fee_currency"
$ws.Cells.Item(1,8).Value = "comment"
Set-HeaderStyle $ws "G1" "H1"
$ws.Columns.Item(8).ColumnWidth = 20.75

$ws.Cells.Item(2,5).Value = "USD"
$ws.Cells.Item(2,6).Value = 0
$ws.Cells.Item(2,7).Value = "USD"

$ws.Cells.Item(3,5).Value = "USD"
$ws.Cells.Item(3,6).Value = 0
$ws.Cells.Item(3,7).Value = "USD"

# ---------------------------------------------------------------------------
# money_transfers: date, buy_date, amount, fees, currency, comment
#                -> date, buy_date, amount, currency, fees, fee_currency, comment
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("money_transfers")

$ws.Cells.Item(1,4).Value = "currency"
$ws.Cells.Item(1,5).Value = "fees"
$ws.Cells.Item(1,6).Value = "fee_currency"
$ws.Cells.Item(1,7).Value = "comment"
Set-HeaderStyle $ws "E1" "F1"
Set-HeaderStyle $ws "D1" "G1"

$ws.Cells.Item(2,4).Value = "USD"
$ws.Cells.Item(2,5).Value = 0
$ws.Cells.Item(2,6).Value = "USD"
Set-HeaderStyle $ws "E2" "F2"

# ---------------------------------------------------------------------------
# sell_orders: date, symbol, quantity, sell_price, fees, currency, comment
#           -> date, symbol, quantity, sell_price, currency, fees, fee_currency, comment
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("sell_orders")

$ws.Cells.Item(1,5).Value = "currency"
$ws.Cells.Item(1,6).Value = "fees"
$ws.Cells.Item(1,7).Value = "fee_currency"
$ws.Cells.Item(1,8).Value = "comment"
Set-HeaderStyle $ws "G1" "H1"
$ws.Columns.Item(8).ColumnWidth = 14.5

$ws.Cells.Item(2,5).Value = "USD"
$ws.Cells.Item(2,6).Value = 0
$ws.Cells.Item(2,7).Value = "USD"

# ---------------------------------------------------------------------------
# currency_conversions: date, source_amount, source_fees, source_currency,
#                        target_amount, target_fees, target_currency, comment
#                     -> date, source_amount, source_currency, target_amount,
#                        target_currency, fees, fee_currency, comment
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("currency_conversions")

$ws.Cells.Item(1,3).Value = "source_currency"
$ws.Cells.Item(1,4).Value = "target_amount"
$ws.Cells.Item(1,5).Value = "target_currency"
$ws.Cells.Item(1,6).Value = "fees"
$ws.Cells.Item(1,7).Value = "fee_currency"

# ---------------------------------------------------------------------------
# Move the active tab from currency_conversions to espp
# ---------------------------------------------------------------------------
$wb.Worksheets.Item("espp").Activate()
